$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-6 from serial 45221 to 45224
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45224
}
